$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: headers
$ws.Range("A1").Value2 = "Row"
$ws.Range("B1").Value2 = "Prognose"
$ws.Range("C1").Value2 = "surveys"
$ws.Range("D1").Value2 = "production"
$ws.Range("E1").Value2 = "orders"
$ws.Range("F1").Value2 = "turnover"
$ws.Range("G1").Value2 = "financial"
$ws.Range("H1").Value2 = "labor market"
$ws.Range("I1").Value2 = "prices"
$ws.Range("J1").Value2 = "national accounts"
$ws.Range("K1").Value2 = "Revision"

# Row 2
$ws.Range("A2").Value2 = "'2025-03-30"
$ws.Range("B2").Value2 = 0.37494215285258325
$ws.Range("C2").Value2 = 0
$ws.Range("D2").Value2 = 0
$ws.Range("E2").Value2 = 0
$ws.Range("F2").Value2 = 0
$ws.Range("G2").Value2 = 0
$ws.Range("H2").Value2 = 0
$ws.Range("I2").Value2 = 0
$ws.Range("J2").Value2 = 0
$ws.Range("K2").Value2 = 0

# Row 3
$ws.Range("A3").Value2 = "'2025-04-15"
$ws.Range("B3").Value2 = 0.31728777927991436
$ws.Range("C3").Value2 = 0
$ws.Range("D3").Value2 = -0.048250834437548649
$ws.Range("E3").Value2 = -0.011237472636205777
$ws.Range("F3").Value2 = 0.0011219728675306727
$ws.Range("G3").Value2 = 0.0050848182041610922
$ws.Range("H3").Value2 = 0.00060397080384812164
$ws.Range("I3").Value2 = -0.010495629456459092
$ws.Range("J3").Value2 = 0
$ws.Range("K3").Value2 = 0.0055188010820047539

# Row 4
$ws.Range("A4").Value2 = "'2025-04-30"
$ws.Range("B4").Value2 = 0.28186685752609764
$ws.Range("C4").Value2 = 0.011297879567615798
$ws.Range("D4").Value2 = 0
$ws.Range("E4").Value2 = -0.0015583423365256245
$ws.Range("F4").Value2 = 0.00051196392761580051
$ws.Range("G4").Value2 = 0
$ws.Range("H4").Value2 = 0.0039928363574255644
$ws.Range("I4").Value2 = -0.046894796157765192
$ws.Range("J4").Value2 = -0.0048044257985845407
$ws.Range("K4").Value2 = 0.0020339626864014915

# Row 5
$ws.Range("A5").Value2 = "'2025-05-15"
$ws.Range("B5").Value2 = 0.4371086494903183
$ws.Range("C5").Value2 = 0.2358094723210559
$ws.Range("D5").Value2 = -0.031193923573388738
$ws.Range("E5").Value2 = -0.0074990983242485612
$ws.Range("F5").Value2 = -0.019366763523345647
$ws.Range("G5").Value2 = -0.010496619648665126
$ws.Range("H5").Value2 = 0.0022815750952709427
$ws.Range("I5").Value2 = -0.018991798149531766
$ws.Range("J5").Value2 = 0
$ws.Range("K5").Value2 = 0.0046989477670736912

# Row 6
$ws.Range("A6").Value2 = "'2025-05-30"
$ws.Range("B6").Value2 = 0.52814190741875011
$ws.Range("C6").Value2 = 0.03926943925387117
$ws.Range("D6").Value2 = 0
$ws.Range("E6").Value2 = -0.0064205297882995883
$ws.Range("F6").Value2 = 0.011350951816908976
$ws.Range("G6").Value2 = 0
$ws.Range("H6").Value2 = 0.0043916781529059331
$ws.Range("I6").Value2 = 0.037373047067549917
$ws.Range("J6").Value2 = 0
$ws.Range("K6").Value2 = 0.0050686714254953857

# Row 7
$ws.Range("A7").Value2 = "'2025-06-15"
$ws.Range("B7").Value2 = 0.54927984900814897
$ws.Range("C7").Value2 = 0
$ws.Range("D7").Value2 = 0.050935570642089791
$ws.Range("E7").Value2 = -0.0029573773074936774
$ws.Range("F7").Value2 = -0.064338452021497641
$ws.Range("G7").Value2 = 0.008063204740199904
$ws.Range("H7").Value2 = 0
$ws.Range("I7").Value2 = -0.0013837665298571545
$ws.Range("J7").Value2 = 0
$ws.Range("K7").Value2 = 0.030818762065957683

# Row 8
$ws.Range("A8").Value2 = "'2025-06-30"
$ws.Range("B8").Value2 = 0.075769554803508943
$ws.Range("C8").Value2 = -0.43937466354771604
$ws.Range("D8").Value2 = 0
$ws.Range("E8").Value2 = 0.0059380149059734821
$ws.Range("F8").Value2 = 0.033003659045883628
$ws.Range("G8").Value2 = 0
$ws.Range("H8").Value2 = -0.001122915937623365
$ws.Range("I8").Value2 = -0.070782846367150251
$ws.Range("J8").Value2 = 0
$ws.Range("K8").Value2 = -0.0011715423040075112

# Row 9
$ws.Range("A9").Value2 = "'2025-07-15"
$ws.Range("B9").Value2 = -0.00056392017379197634
$ws.Range("C9").Value2 = 0
$ws.Range("D9").Value2 = -0.10473246181758607
$ws.Range("E9").Value2 = -0.066448226700680901
$ws.Range("F9").Value2 = 0.092894759177118494
$ws.Range("G9").Value2 = -0.0039845070481038906
$ws.Range("H9").Value2 = [double]"-7.6154081641479753e-05"
$ws.Range("I9").Value2 = -0.00071605479960878374
$ws.Range("J9").Value2 = 0
$ws.Range("K9").Value2 = 0.0067291702932017333

# Row 10
$ws.Range("A10").Value2 = "'2025-07-30"
$ws.Range("B10").Value2 = 0.25447748280026017
$ws.Range("C10").Value2 = 0.26341824452223311
$ws.Range("D10").Value2 = 0
$ws.Range("E10").Value2 = 0.0098352589657540943
$ws.Range("F10").Value2 = -0.0057059535851061742
$ws.Range("G10").Value2 = 0
$ws.Range("H10").Value2 = -0.0032438949597869003
$ws.Range("I10").Value2 = 0.011991957638635673
$ws.Range("J10").Value2 = -0.010737831459777646
$ws.Range("K10").Value2 = -0.010516378147900074

# Row 11
$ws.Range("A11").Value2 = "'2025-08-15"
$ws.Range("B11").Value2 = 0.6012162683384028
$ws.Range("C11").Value2 = 0
$ws.Range("D11").Value2 = 0.27608880179844986
$ws.Range("E11").Value2 = 0.0086376429710910314
$ws.Range("F11").Value2 = 0.018996936426581873
$ws.Range("G11").Value2 = -0.0055586309381815153
$ws.Range("H11").Value2 = -0.0010836579988595192
$ws.Range("I11").Value2 = 0.087215307302099823
$ws.Range("J11").Value2 = 0
$ws.Range("K11").Value2 = -0.037557614023038921

# Row 12
$ws.Range("A12").Value2 = "'2025-08-30"
$ws.Range("B12").Value2 = 0.23692507859600248
$ws.Range("C12").Value2 = -0.32331421480764799
$ws.Range("D12").Value2 = 0
$ws.Range("E12").Value2 = 0.0033425990038162671
$ws.Range("F12").Value2 = 0.00011646454896633338
$ws.Range("G12").Value2 = 0
$ws.Range("H12").Value2 = [double]"-7.2634137113761788e-05"
$ws.Range("I12").Value2 = -0.027620400336508352
$ws.Range("J12").Value2 = 0
$ws.Range("K12").Value2 = -0.016743004013912821
